# The commit swaps the presentation's theme from the "Integral" palette to
# the stock "Office Theme" palette (a Design-tab theme change). PowerPoint
# persists this by writing the new theme colours into the theme part that
# the slide master references (ppt/theme/theme1.xml).
#
# ThemeColorScheme items are ordered: dk1, lt1, dk2, lt2, accent1..accent6,
# hlink, folHlink - matching <a:clrScheme> child order. RGB is the packed
# 0xBBGGRR long (same convention as the VBA RGB() function).

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

# Office Theme colour scheme (target palette)
$colors.Item(1).RGB  = 0          # dk1      000000
$colors.Item(2).RGB  = 16777215   # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388    # dk2      44546A
$colors.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501    # accent2  ED7D31
$colors.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colors.Item(8).RGB  = 49407      # accent4  FFC000
$colors.Item(9).RGB  = 12874308   # accent5  4472C4
$colors.Item(10).RGB = 4697456    # accent6  70AD47
$colors.Item(11).RGB = 12673797   # hlink    0563C1
$colors.Item(12).RGB = 7491477    # folHlink 954F72
